$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume/Hora columns for data rows 2-51 keep their original
# text (string) storage instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Cells.Item(2,4).Value = '306.86'
$ws.Cells.Item(2,5).Value = '2.37%'
$ws.Cells.Item(2,7).Value = '15'
$ws.Cells.Item(3,4).Value = '35.88'
$ws.Cells.Item(3,5).Value = '1.45%'
$ws.Cells.Item(3,7).Value = '15'
$ws.Cells.Item(4,4).Value = '5.080'
$ws.Cells.Item(4,5).Value = '0.87%'
$ws.Cells.Item(4,7).Value = '15'
$ws.Cells.Item(5,4).Value = '0.08087'
$ws.Cells.Item(5,5).Value = '1.74%'
$ws.Cells.Item(5,7).Value = '15'
$ws.Cells.Item(6,4).Value = '1.935'
$ws.Cells.Item(6,5).Value = '2.11%'
$ws.Cells.Item(6,7).Value = '15'
$ws.Cells.Item(7,4).Value = '4.169'
$ws.Cells.Item(7,5).Value = '2.91%'
$ws.Cells.Item(7,7).Value = '15'
$ws.Cells.Item(8,4).Value = '7.826'
$ws.Cells.Item(8,5).Value = '0.55%'
$ws.Cells.Item(8,7).Value = '15'
$ws.Cells.Item(9,4).Value = '0.9349'
$ws.Cells.Item(9,5).Value = '1.39%'
$ws.Cells.Item(9,7).Value = '15'
$ws.Cells.Item(10,4).Value = '0.1333'
$ws.Cells.Item(10,5).Value = '-9.28%'
$ws.Cells.Item(10,7).Value = '15'
$ws.Cells.Item(11,4).Value = '0.1909'
$ws.Cells.Item(11,5).Value = '1.01%'
$ws.Cells.Item(11,7).Value = '15'
$ws.Cells.Item(12,4).Value = '0.09205'
$ws.Cells.Item(12,5).Value = '0.88%'
$ws.Cells.Item(12,7).Value = '15'
$ws.Cells.Item(13,4).Value = '0.03523'
$ws.Cells.Item(13,5).Value = '2.49%'
$ws.Cells.Item(13,7).Value = '15'
$ws.Cells.Item(14,4).Value = '0.09892'
$ws.Cells.Item(14,5).Value = '0.13%'
$ws.Cells.Item(14,7).Value = '15'
$ws.Cells.Item(15,4).Value = '0.001431'
$ws.Cells.Item(15,5).Value = '2.95%'
$ws.Cells.Item(15,7).Value = '15'
$ws.Cells.Item(16,4).Value = '0.005874'
$ws.Cells.Item(16,5).Value = '0.62%'
$ws.Cells.Item(16,7).Value = '15'
$ws.Cells.Item(17,4).Value = '3.604'
$ws.Cells.Item(17,5).Value = '2.55%'
$ws.Cells.Item(17,7).Value = '15'
$ws.Cells.Item(18,7).Value = '15'
$ws.Cells.Item(19,4).Value = '0.3455'
$ws.Cells.Item(19,5).Value = '1.48%'
$ws.Cells.Item(19,7).Value = '15'
$ws.Cells.Item(20,4).Value = '0.1346'
$ws.Cells.Item(20,5).Value = '3.87%'
$ws.Cells.Item(20,7).Value = '15'
$ws.Cells.Item(21,4).Value = '5.199'
$ws.Cells.Item(21,5).Value = '2.46%'
$ws.Cells.Item(21,7).Value = '15'
$ws.Cells.Item(22,4).Value = '0.2627'
$ws.Cells.Item(22,5).Value = '9.32%'
$ws.Cells.Item(22,7).Value = '15'
$ws.Cells.Item(23,4).Value = '0.04393'
$ws.Cells.Item(23,5).Value = '-1.49%'
$ws.Cells.Item(23,7).Value = '15'
$ws.Cells.Item(24,4).Value = '0.001242'
$ws.Cells.Item(24,5).Value = '1.97%'
$ws.Cells.Item(24,7).Value = '15'
$ws.Cells.Item(25,4).Value = '0.004754'
$ws.Cells.Item(25,5).Value = '-0.12%'
$ws.Cells.Item(25,7).Value = '15'
$ws.Cells.Item(26,4).Value = '0.0001305'
$ws.Cells.Item(26,5).Value = '5.68%'
$ws.Cells.Item(26,7).Value = '15'
$ws.Cells.Item(27,5).Value = '4.26%'
$ws.Cells.Item(27,7).Value = '15'
$ws.Cells.Item(28,7).Value = '15'
$ws.Cells.Item(29,7).Value = '15'
$ws.Cells.Item(30,7).Value = '15'
$ws.Cells.Item(31,7).Value = '15'
$ws.Cells.Item(32,7).Value = '15'
$ws.Cells.Item(33,7).Value = '15'
$ws.Cells.Item(34,7).Value = '15'
$ws.Cells.Item(35,7).Value = '15'
$ws.Cells.Item(36,7).Value = '15'
$ws.Cells.Item(37,7).Value = '15'
$ws.Cells.Item(38,7).Value = '15'
$ws.Cells.Item(39,4).Value = '0.01988'
$ws.Cells.Item(39,5).Value = '4.52%'
$ws.Cells.Item(39,7).Value = '15'
$ws.Cells.Item(40,4).Value = '0.05002'
$ws.Cells.Item(40,7).Value = '15'
$ws.Cells.Item(41,4).Value = '0.01124'
$ws.Cells.Item(41,5).Value = '15.43%'
$ws.Cells.Item(41,7).Value = '15'
$ws.Cells.Item(42,4).Value = '0.007644'
$ws.Cells.Item(42,5).Value = '3.63%'
$ws.Cells.Item(42,7).Value = '15'
$ws.Cells.Item(43,4).Value = '0.1376'
$ws.Cells.Item(43,5).Value = '3.84%'
$ws.Cells.Item(43,7).Value = '15'
$ws.Cells.Item(44,4).Value = '0.002109'
$ws.Cells.Item(44,5).Value = '-0.47%'
$ws.Cells.Item(44,7).Value = '15'
$ws.Cells.Item(45,4).Value = '0.01136'
$ws.Cells.Item(45,5).Value = '21.55%'
$ws.Cells.Item(45,7).Value = '15'
$ws.Cells.Item(46,4).Value = '0.00006422'
$ws.Cells.Item(46,5).Value = '2.23%'
$ws.Cells.Item(46,7).Value = '15'
$ws.Cells.Item(47,4).Value = '0.00000000753'
$ws.Cells.Item(47,5).Value = '0.01%'
$ws.Cells.Item(47,7).Value = '15'
$ws.Cells.Item(48,4).Value = '65.22'
$ws.Cells.Item(48,5).Value = '1.15%'
$ws.Cells.Item(48,7).Value = '15'
$ws.Cells.Item(49,5).Value = '-28.25%'
$ws.Cells.Item(49,7).Value = '15'
$ws.Cells.Item(50,4).Value = '0.00002109'
$ws.Cells.Item(50,5).Value = '0.01%'
$ws.Cells.Item(50,7).Value = '15'
$ws.Cells.Item(51,4).Value = '0.0002009'
$ws.Cells.Item(51,5).Value = '0.01%'
$ws.Cells.Item(51,7).Value = '15'
